$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 68. This shifts the existing
# rows 68-175 down to 69-176 (and the sheet dimension grows to R176),
# exactly like pressing "Insert" on the row header in Excel.
$ws.Rows(68).Insert()

# Populate the newly inserted row 68 with the new weekly record. All
# columns besides the date (D) and volume (J) keep the same values the
# old row 68 had (Terminal Hortofrutícola Agro Chillán / Ñuble / Repollo
# / Crespo record / Primera / 600-700, avg 650 / Provincia de Diguillín).
$ws.Cells.Item(68, 1).Value = 7
$ws.Cells.Item(68, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(68, 3).Value = "Ñuble"
$ws.Cells.Item(68, 4).Value = 44580
$ws.Cells.Item(68, 5).Value = 16
$ws.Cells.Item(68, 6).Value = 100112006
$ws.Cells.Item(68, 7).Value = "Repollo"
$ws.Cells.Item(68, 8).Value = "Crespo record"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 240
$ws.Cells.Item(68, 11).Value = 600
$ws.Cells.Item(68, 12).Value = 700
$ws.Cells.Item(68, 13).Value = 650
$ws.Cells.Item(68, 14).Value = "`$/unidad"
$ws.Cells.Item(68, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(68, 16).Value = 650
$ws.Cells.Item(68, 17).Value = 1
$ws.Cells.Item(68, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# rest of column D.
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(69, 4).NumberFormat
